$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string "Resolving-Mac" (index 25) is introduced implicitly by
# assigning it as a cell value below (D5 is the first cell to reference it).

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cntn2"
$ws.Cells.Item(2, 3).Value = "Nrcam"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2862826666666667
$ws.Cells.Item(2, 8).Value = 0.8588480000000001
$ws.Cells.Item(2, 9).Value = 0.560705294934871
$ws.Cells.Item(2, 10).Value = 0.560705294934871
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.008656666666666667
$ws.Cells.Item(2, 14).Value = 0.02597
$ws.Cells.Item(2, 15).Value = 0.006895973038524511
$ws.Cells.Item(2, 16).Value = 0.006895973038524512
$ws.Cells.Item(2, 17).Value = 0.002478253617777778
$ws.Cells.Item(2, 18).Value = 0.02230428256
$ws.Cells.Item(2, 19).Value = 0.003866608596428804
$ws.Cells.Item(2, 20).Value = 0.003866608596428805

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cntn2"
$ws.Cells.Item(3, 3).Value = "Nrcam"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2862826666666667
$ws.Cells.Item(3, 8).Value = 0.8588480000000001
$ws.Cells.Item(3, 9).Value = 0.560705294934871
$ws.Cells.Item(3, 10).Value = 0.560705294934871
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.1287113333333333
$ws.Cells.Item(3, 14).Value = 0.386134
$ws.Cells.Item(3, 15).Value = 0.1025325241916682
$ws.Cells.Item(3, 16).Value = 0.1025325241916682
$ws.Cells.Item(3, 17).Value = 0.03684782373688889
$ws.Cells.Item(3, 18).Value = 0.331630413632
$ws.Cells.Item(3, 19).Value = 0.05749052921730612
$ws.Cells.Item(3, 20).Value = 0.05749052921730612

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Cntn2"
$ws.Cells.Item(4, 3).Value = "Nrcam"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.2862826666666667
$ws.Cells.Item(4, 8).Value = 0.8588480000000001
$ws.Cells.Item(4, 9).Value = 0.560705294934871
$ws.Cells.Item(4, 10).Value = 0.560705294934871
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.089774666666667
$ws.Cells.Item(4, 14).Value = 3.269324
$ws.Cells.Item(4, 15).Value = 0.8681236102503316
$ws.Cells.Item(4, 16).Value = 0.8681236102503316
$ws.Cells.Item(4, 17).Value = 0.3119835976391112
$ws.Cells.Item(4, 18).Value = 2.807852378752
$ws.Cells.Item(4, 19).Value = 0.4867615049253372
$ws.Cells.Item(4, 20).Value = 0.4867615049253372

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Cntn2"
$ws.Cells.Item(5, 3).Value = "Nrcam"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.2862826666666667
$ws.Cells.Item(5, 8).Value = 0.8588480000000001
$ws.Cells.Item(5, 9).Value = 0.560705294934871
$ws.Cells.Item(5, 10).Value = 0.560705294934871
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.02817933333333333
$ws.Cells.Item(5, 14).Value = 0.084538
$ws.Cells.Item(5, 15).Value = 0.02244789251947575
$ws.Cells.Item(5, 16).Value = 0.02244789251947575
$ws.Cells.Item(5, 17).Value = 0.008067254691555557
$ws.Cells.Item(5, 18).Value = 0.072605292224
$ws.Cells.Item(5, 19).Value = 0.01258665219579893
$ws.Cells.Item(5, 20).Value = 0.01258665219579893

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Cntn2"
$ws.Cells.Item(6, 3).Value = "Nrcam"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.1092446666666667
$ws.Cells.Item(6, 8).Value = 0.327734
$ws.Cells.Item(6, 9).Value = 0.2139635757784672
$ws.Cells.Item(6, 10).Value = 0.2139635757784672
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.008656666666666667
$ws.Cells.Item(6, 14).Value = 0.02597
$ws.Cells.Item(6, 15).Value = 0.006895973038524511
$ws.Cells.Item(6, 16).Value = 0.006895973038524512
$ws.Cells.Item(6, 17).Value = 0.0009456946644444445
$ws.Cells.Item(6, 18).Value = 0.00851125198
$ws.Cells.Item(6, 19).Value = 0.001475487049794606
$ws.Cells.Item(6, 20).Value = 0.001475487049794606

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Cntn2"
$ws.Cells.Item(7, 3).Value = "Nrcam"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.1092446666666667
$ws.Cells.Item(7, 8).Value = 0.327734
$ws.Cells.Item(7, 9).Value = 0.2139635757784672
$ws.Cells.Item(7, 10).Value = 0.2139635757784672
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.1287113333333333
$ws.Cells.Item(7, 14).Value = 0.386134
$ws.Cells.Item(7, 15).Value = 0.1025325241916682
$ws.Cells.Item(7, 16).Value = 0.1025325241916682
$ws.Cells.Item(7, 17).Value = 0.01406102670622222
$ws.Cells.Item(7, 18).Value = 0.126549240356
$ws.Cells.Item(7, 19).Value = 0.02193822550964152
$ws.Cells.Item(7, 20).Value = 0.02193822550964153

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Cntn2"
$ws.Cells.Item(8, 3).Value = "Nrcam"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.1092446666666667
$ws.Cells.Item(8, 8).Value = 0.327734
$ws.Cells.Item(8, 9).Value = 0.2139635757784672
$ws.Cells.Item(8, 10).Value = 0.2139635757784672
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.089774666666667
$ws.Cells.Item(8, 14).Value = 3.269324
$ws.Cells.Item(8, 15).Value = 0.8681236102503316
$ws.Cells.Item(8, 16).Value = 0.8681236102503316
$ws.Cells.Item(8, 17).Value = 0.1190520702017778
$ws.Cells.Item(8, 18).Value = 1.071468631816
$ws.Cells.Item(8, 19).Value = 0.1857468318668734
$ws.Cells.Item(8, 20).Value = 0.1857468318668734

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Cntn2"
$ws.Cells.Item(9, 3).Value = "Nrcam"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.1092446666666667
$ws.Cells.Item(9, 8).Value = 0.327734
$ws.Cells.Item(9, 9).Value = 0.2139635757784672
$ws.Cells.Item(9, 10).Value = 0.2139635757784672
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.02817933333333333
$ws.Cells.Item(9, 14).Value = 0.084538
$ws.Cells.Item(9, 15).Value = 0.02244789251947575
$ws.Cells.Item(9, 16).Value = 0.02244789251947575
$ws.Cells.Item(9, 17).Value = 0.003078441876888889
$ws.Cells.Item(9, 18).Value = 0.027705976892
$ws.Cells.Item(9, 19).Value = 0.004803031352157736
$ws.Cells.Item(9, 20).Value = 0.004803031352157736

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Cntn2"
$ws.Cells.Item(10, 3).Value = "Nrcam"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.1150486666666666
$ws.Cells.Item(10, 8).Value = 0.345146
$ws.Cells.Item(10, 9).Value = 0.2253311292866618
$ws.Cells.Item(10, 10).Value = 0.2253311292866618
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.008656666666666667
$ws.Cells.Item(10, 14).Value = 0.02597
$ws.Cells.Item(10, 15).Value = 0.006895973038524511
$ws.Cells.Item(10, 16).Value = 0.006895973038524512
$ws.Cells.Item(10, 17).Value = 0.0009959379577777776
$ws.Cells.Item(10, 18).Value = 0.008963441619999998
$ws.Cells.Item(10, 19).Value = 0.001553877392301101
$ws.Cells.Item(10, 20).Value = 0.001553877392301101

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Cntn2"
$ws.Cells.Item(11, 3).Value = "Nrcam"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.1150486666666666
$ws.Cells.Item(11, 8).Value = 0.345146
$ws.Cells.Item(11, 9).Value = 0.2253311292866618
$ws.Cells.Item(11, 10).Value = 0.2253311292866618
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.1287113333333333
$ws.Cells.Item(11, 14).Value = 0.386134
$ws.Cells.Item(11, 15).Value = 0.1025325241916682
$ws.Cells.Item(11, 16).Value = 0.1025325241916682
$ws.Cells.Item(11, 17).Value = 0.01480806728488888
$ws.Cells.Item(11, 18).Value = 0.133272605564
$ws.Cells.Item(11, 19).Value = 0.02310376946472057
$ws.Cells.Item(11, 20).Value = 0.02310376946472057

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Cntn2"
$ws.Cells.Item(12, 3).Value = "Nrcam"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.1150486666666666
$ws.Cells.Item(12, 8).Value = 0.345146
$ws.Cells.Item(12, 9).Value = 0.2253311292866618
$ws.Cells.Item(12, 10).Value = 0.2253311292866618
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.089774666666667
$ws.Cells.Item(12, 14).Value = 3.269324
$ws.Cells.Item(12, 15).Value = 0.8681236102503316
$ws.Cells.Item(12, 16).Value = 0.8681236102503316
$ws.Cells.Item(12, 17).Value = 0.1253771223671111
$ws.Cells.Item(12, 18).Value = 1.128394101304
$ws.Cells.Item(12, 19).Value = 0.195615273458121
$ws.Cells.Item(12, 20).Value = 0.1956152734581211

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Cntn2"
$ws.Cells.Item(13, 3).Value = "Nrcam"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.1150486666666666
$ws.Cells.Item(13, 8).Value = 0.345146
$ws.Cells.Item(13, 9).Value = 0.2253311292866618
$ws.Cells.Item(13, 10).Value = 0.2253311292866618
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.02817933333333333
$ws.Cells.Item(13, 14).Value = 0.084538
$ws.Cells.Item(13, 15).Value = 0.02244789251947575
$ws.Cells.Item(13, 16).Value = 0.02244789251947575
$ws.Cells.Item(13, 17).Value = 0.003241994727555555
$ws.Cells.Item(13, 18).Value = 0.029177952548
$ws.Cells.Item(13, 19).Value = 0.005058208971519078
$ws.Cells.Item(13, 20).Value = 0.005058208971519079

Write-Output "Updated rows 2-13 with new TPM values and added Resolving-Mac cluster"
